# Apply the edits described by the diff:
#  1. "Presentado a la instructora:" -> "Presentado al instructor:"
#  2. "Elizabeth Robayo Ramirez" -> "Gustavo Adolfo Rodriguez"
#     (and the _GoBack bookmark now wraps this new name, moving away
#      from its old location further down the document)
#  3. Date "10 de Diciembre del 2023" -> "19 de Febrero del 2024"

$d = $word.ActiveDocument

# 1. Fix the greeting line: "instructora" (feminine) -> "instructor" (masculine)
$d.Content.Find.Execute(
    "Presentado a la instructora:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Presentado al instructor:", 2) | Out-Null

# 2. Replace the previous instructor's name with the new one
$d.Content.Find.Execute(
    "Elizabeth Robayo Ramirez", $true, $false, $false, $false, $false,
    $true, 1, $false, "Gustavo Adolfo Rodriguez", 2) | Out-Null

# 3. Update the submission date
$d.Content.Find.Execute(
    "10 de Diciembre del 2023", $true, $false, $false, $false, $false,
    $true, 1, $false, "19 de Febrero del 2024", 2) | Out-Null

# 4. Move the "_GoBack" bookmark (Word's last-edit marker) onto the
#    paragraph that now holds the new instructor's name - this also
#    removes it from its previous location further down the document,
#    exactly as Word does when a single-named bookmark is re-added.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "Gustavo Adolfo Rodriguez") {
        $nameRange = $p.Range
        $nameRange.MoveEnd(1, -1) | Out-Null
        $d.Bookmarks.Add("_GoBack", $nameRange) | Out-Null
        break
    }
}
